# Update Wnt9a-Fzd10 LR-pairs sheet with newly recomputed TPM figures.
# The refreshed pipeline now also reports Target cluster = "ECs" rows for
# each sending cluster (previously only "MuSCs" targets were emitted), so
# the 3-row table grows to 6 rows and most numeric columns are recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending=ECs, Target=ECs (was Target=MuSCs; values recomputed)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt9a"
$ws.Range("C2").Value = "Fzd10"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2304126666666667
$ws.Range("H2").Value = 0.691238
$ws.Range("I2").Value = 0.03265479005310033
$ws.Range("J2").Value = 0.03265479005310033
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.003092666666666667
$ws.Range("N2").Value = 0.009278
$ws.Range("O2").Value = 0.03934390080485799
$ws.Range("P2").Value = 0.03934390080485798
$ws.Range("Q2").Value = 0.0007125895737777778
$ws.Range("R2").Value = 0.006413306164
$ws.Range("S2").Value = 0.001284766820652643
$ws.Range("T2").Value = 0.001284766820652642

# Row 3: Sending=ECs, Target=MuSCs (was Sending=FAPs; values recomputed)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt9a"
$ws.Range("C3").Value = "Fzd10"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2304126666666667
$ws.Range("H3").Value = 0.691238
$ws.Range("I3").Value = 0.03265479005310033
$ws.Range("J3").Value = 0.03265479005310033
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.07551333333333334
$ws.Range("N3").Value = 0.22654
$ws.Range("O3").Value = 0.9606560991951421
$ws.Range("P3").Value = 0.9606560991951419
$ws.Range("Q3").Value = 0.01739922850222222
$ws.Range("R3").Value = 0.15659305652
$ws.Range("S3").Value = 0.03137002323244768
$ws.Range("T3").Value = 0.03137002323244768

# Row 4: Sending=FAPs, Target=ECs (was Sending=MuSCs; values recomputed)
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt9a"
$ws.Range("C4").Value = "Fzd10"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.016075000000001
$ws.Range("H4").Value = 18.048225
$ws.Range("I4").Value = 0.8526166070240881
$ws.Range("J4").Value = 0.8526166070240883
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.003092666666666667
$ws.Range("N4").Value = 0.009278
$ws.Range("O4").Value = 0.03934390080485799
$ws.Range("P4").Value = 0.03934390080485798
$ws.Range("Q4").Value = 0.01860571461666667
$ws.Range("R4").Value = 0.16745143155
$ws.Range("S4").Value = 0.03354526321133031
$ws.Range("T4").Value = 0.03354526321133031

# Row 5 (new): Sending=FAPs, Target=MuSCs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt9a"
$ws.Range("C5").Value = "Fzd10"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.016075000000001
$ws.Range("H5").Value = 18.048225
$ws.Range("I5").Value = 0.8526166070240881
$ws.Range("J5").Value = 0.8526166070240883
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.07551333333333334
$ws.Range("N5").Value = 0.22654
$ws.Range("O5").Value = 0.9606560991951421
$ws.Range("P5").Value = 0.9606560991951419
$ws.Range("Q5").Value = 0.4542938768333334
$ws.Range("R5").Value = 4.0886448915
$ws.Range("S5").Value = 0.8190713438127579
$ws.Range("T5").Value = 0.8190713438127579

# Row 6 (new): Sending=MuSCs, Target=ECs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Wnt9a"
$ws.Range("C6").Value = "Fzd10"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.8095266666666667
$ws.Range("H6").Value = 2.42858
$ws.Range("I6").Value = 0.1147286029228115
$ws.Range("J6").Value = 0.1147286029228115
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.003092666666666667
$ws.Range("N6").Value = 0.009278
$ws.Range("O6").Value = 0.03934390080485799
$ws.Range("P6").Value = 0.03934390080485798
$ws.Range("Q6").Value = 0.002503596137777778
$ws.Range("R6").Value = 0.02253236524
$ws.Range("S6").Value = 0.004513870772875036
$ws.Range("T6").Value = 0.004513870772875035

# Row 7 (new): Sending=MuSCs, Target=MuSCs (was row 4's content before the shift)
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Wnt9a"
$ws.Range("C7").Value = "Fzd10"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.8095266666666667
$ws.Range("H7").Value = 2.42858
$ws.Range("I7").Value = 0.1147286029228115
$ws.Range("J7").Value = 0.1147286029228115
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.07551333333333334
$ws.Range("N7").Value = 0.22654
$ws.Range("O7").Value = 0.9606560991951421
$ws.Range("P7").Value = 0.9606560991951419
$ws.Range("Q7").Value = 0.06113005702222223
$ws.Range("R7").Value = 0.5501705132
$ws.Range("S7").Value = 0.1102147321499365
$ws.Range("T7").Value = 0.1102147321499365
